# "break out stock.yaml completed"
#
# - Sheet "10per change": append row 21 (APLAPOLLO, 13/06/2024 08:45:34) and
#   convert the existing E20 "bsecode" cell from text to a real number.
# - Sheet "3 V 0.3": append row 7 (PAISALO, 13/06/2024 08:45:34).
#
# Note: the "bsecode" column (E) is stored as literal text in the new rows
# (matching the source data quirk), even though it looks numeric. Assigning
# a numeric-looking string straight to `.Value` gets auto-coerced to a
# number by this host, same as real Excel. To land a genuine text cell
# without leaving a stray NumberFormat="@" style behind, we build the text
# with TEXT() in a scratch cell, Copy it, and PasteSpecial only the value
# into the destination (this carries over the text type but not any
# formatting), then clean up the scratch cell.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("10per change")

$ws1.Range("A21").Value = "13/06/2024 08:45:34"
$ws1.Range("B21").Value = 1
$ws1.Range("C21").Value = "APLAPOLLO"
$ws1.Range("D21").Value = "Apl Apollo Tubes Limited"

$ws1.Range("Z1").Formula = '=TEXT(533758,"0")'
$ws1.Range("Z1").Copy()
$ws1.Range("E21").PasteSpecial(-4163)
$ws1.Range("Z1").Clear()

$ws1.Range("F21").Value = -2.43
$ws1.Range("G21").Value = 1542.4
$ws1.Range("H21").Value = 466997

# E20 was text "533758" in the source; the edit turns it into a real number.
$ws1.Range("E20").Value = 533758

$ws2 = $wb.Worksheets.Item("3 V 0.3")

$ws2.Range("A7").Value = "13/06/2024 08:45:34"
$ws2.Range("B7").Value = 1
$ws2.Range("C7").Value = "PAISALO"
$ws2.Range("D7").Value = "Paisalo Digital Ltd"

$ws2.Range("Z1").Formula = '=TEXT(532900,"0")'
$ws2.Range("Z1").Copy()
$ws2.Range("E7").PasteSpecial(-4163)
$ws2.Range("Z1").Clear()

$ws2.Range("F7").Value = 8.03
$ws2.Range("G7").Value = 74.52
$ws2.Range("H7").Value = 4339439
